$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "20.391.91"
$ws.Range("E2").Value = "  +2.25%  "

# Row 3
$ws.Range("D3").Value = "1.461.75"
$ws.Range("E3").Value = "  +3.40%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.83%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9495"
$ws.Range("E5").Value = "  -5.06%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "274.78"
$ws.Range("E6").Value = "  -0.41%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3650"
$ws.Range("E7").Value = "  -0.45%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3075"
$ws.Range("E8").Value = "  -0.91%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.71"
$ws.Range("E9").Value = "  -0.02%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.037"
$ws.Range("E10").Value = "  +0.08%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06577"
$ws.Range("E11").Value = "  +0.80%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.02%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.424"
$ws.Range("E13").Value = "  -1.26%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.84"
$ws.Range("E14").Value = "  +1.36%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.126"
$ws.Range("E15").Value = "  -1.11%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001025"
$ws.Range("E16").Value = "  +0.73%  "

# Row 17
$ws.Range("D17").Value = "1.463.22"
$ws.Range("E17").Value = "  +3.68%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9678"
$ws.Range("E18").Value = "  -3.26%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05792"
$ws.Range("E19").Value = "  +2.29%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.68"
$ws.Range("E20").Value = "  -1.95%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.438"
$ws.Range("E21").Value = "  -3.19%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.43"
$ws.Range("E22").Value = "  -1.97%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.88"
$ws.Range("E23").Value = "  -0.39%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.246"
$ws.Range("E24").Value = "  +0.36%  "

# Row 25
$ws.Range("D25").Value = "20.436.43"
$ws.Range("E25").Value = "  +2.34%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.81"
$ws.Range("E26").Value = "  +6.68%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.079"
$ws.Range("E27").Value = "  -7.95%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.15"
$ws.Range("E28").Value = "  -0.90%  "

# Row 29
$ws.Range("D29").Value = "1.616.30"
$ws.Range("E29").Value = "  +3.00%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "111.99"
$ws.Range("E30").Value = "  +1.88%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.828"
$ws.Range("E31").Value = "  -1.53%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.876"
$ws.Range("E32").Value = "  -7.62%  "

# Row 33
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.07865"
$ws.Range("E33").Value = "  +2.33%  "

# Row 34
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7905"
$ws.Range("E34").Value = "  -3.28%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.516"
$ws.Range("E35").Value = "  +2.30%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05699"
$ws.Range("E36").Value = "  -1.33%  "

# Row 37
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.141"
$ws.Range("E37").Value = "  +4.36%  "

# Row 38
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.675"
$ws.Range("E38").Value = "  -5.02%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02029"
$ws.Range("E39").Value = "  -1.43%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9569"
$ws.Range("E40").Value = "  -4.39%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.33"
$ws.Range("E41").Value = "  -1.21%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.523"
$ws.Range("E42").Value = "  -9.77%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1853"
$ws.Range("E43").Value = "  -1.80%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5261"
$ws.Range("E44").Value = "  -1.00%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.487"
$ws.Range("E45").Value = "  -1.41%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.96"
$ws.Range("E46").Value = "  -3.95%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "117.29"
$ws.Range("E47").Value = "  +1.63%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5138"
$ws.Range("E48").Value = "  -1.02%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.750"
$ws.Range("E49").Value = "  -0.99%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06417"
$ws.Range("E50").Value = "  +3.55%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9870"
$ws.Range("E51").Value = "  -1.29%  "
